$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.659.59'
$ws.Range('E2').Value = '  -1.90%  '

$ws.Range('D3').Value = '2.295.13'
$ws.Range('E3').Value = '  -2.99%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = "'539.38"
$ws.Range('E5').Value = '  -1.38%  '

$ws.Range('D6').Value = "'128.20"
$ws.Range('E6').Value = '  -3.19%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').Value = "'0.566"
$ws.Range('E8').Value = '  -3.57%  '

$ws.Range('D9').Value = '2.293.76'
$ws.Range('E9').Value = '  -2.93%  '

$ws.Range('D10').Value = "'0.100"
$ws.Range('E10').Value = '  -1.49%  '

$ws.Range('D11').Value = "'5.49"
$ws.Range('E11').Value = '  -0.32%  '

$ws.Range('E12').Value = '  -0.99%  '

$ws.Range('D13').Value = "'0.329"
$ws.Range('E13').Value = '  -2.02%  '

$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = "'23.06"
$ws.Range('E14').Value = '  -4.34%  '

$ws.Range('D15').Value = '59.602.35'
$ws.Range('E15').Value = '  -1.88%  '

$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.704.19'
$ws.Range('E16').Value = '  -2.99%  '

$ws.Range('D17').Value = "'0.0000131"
$ws.Range('E17').Value = '  -1.94%  '

$ws.Range('D18').Value = '2.310.13'
$ws.Range('E18').Value = '  -2.60%  '

$ws.Range('D19').Value = "'10.36"
$ws.Range('E19').Value = '  -3.53%  '

$ws.Range('D20').Value = "'309.73"
$ws.Range('E20').Value = '  -2.48%  '

$ws.Range('D21').Value = "'3.99"

$ws.Range('D22').Value = "'6.43"
$ws.Range('E22').Value = '  -7.19%  '

$ws.Range('D23').Value = "'1.00"
$ws.Range('E23').Value = '  +0.12%  '

$ws.Range('D24').Value = "'63.40"
$ws.Range('E24').Value = '  -0.24%  '

$ws.Range('D25').Value = "'0.168"
$ws.Range('E25').Value = '  -2.41%  '

$ws.Range('E26').Value = '  -0.25%  '

$ws.Range('D27').Value = "'7.69"
$ws.Range('E27').Value = '  -4.09%  '

$ws.Range('D28').Value = "'1.33"
$ws.Range('E28').Value = '  -2.28%  '

$ws.Range('D29').Value = "'171.30"
$ws.Range('E29').Value = '  -0.31%  '

$ws.Range('D30').Value = "'1.17"
$ws.Range('E30').Value = '  +2.23%  '

$ws.Range('D31').Value = "'1.69"
$ws.Range('E31').Value = '  -2.71%  '

$ws.Range('D32').Value = '0.0₃0709'
$ws.Range('E32').Value = '  -4.12%  '

$ws.Range('D33').Value = "'5.77"
$ws.Range('E33').Value = '  -2.29%  '

$ws.Range('E34').Value = '  -5.15%  '

$ws.Range('D35').Value = "'0.376"
$ws.Range('E35').Value = '  -1.83%  '

$ws.Range('E36').Value = '  -0.01%  '

$ws.Range('D37').Value = "'17.70"
$ws.Range('E37').Value = '  -2.11%  '

$ws.Range('E38').Value = '  +0.13%  '

$ws.Range('D39').Value = "'4.02"
$ws.Range('E39').Value = '  -4.44%  '

$ws.Range('D40').Value = "'308.46"
$ws.Range('E40').Value = '  -3.69%  '

$ws.Range('D41').Value = "'37.80"
$ws.Range('E41').Value = '  -1.29%  '

$ws.Range('D42').Value = "'1.49"
$ws.Range('E42').Value = '  -3.73%  '

$ws.Range('D43').Value = "'135.72"
$ws.Range('E43').Value = '  -5.62%  '

$ws.Range('D44').Value = "'3.39"
$ws.Range('E44').Value = '  -2.65%  '

$ws.Range('D45').Value = "'0.0932"
$ws.Range('E45').Value = '  -2.57%  '

$ws.Range('D46').Value = "'0.563"
$ws.Range('E46').Value = '  -0.26%  '

$ws.Range('D47').Value = "'18.39"
$ws.Range('E47').Value = '  -5.29%  '

$ws.Range('D48').Value = "'0.0485"
$ws.Range('E48').Value = '  -3.23%  '

$ws.Range('D49').Value = '0.0₆0217'
$ws.Range('E49').Value = '  +6.35%  '

$ws.Range('E50').Value = '  -1.80%  '

$ws.Range('D51').Value = "'11.00"
$ws.Range('E51').Value = '  -0.40%  '
